$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("Q2").Value = 2.25
$ws.Range("R2").Value = 1.62
